# Journal de travail — add the "creation MCD / git / MLD" session entries
# (row 5 follow-up + brand-new row 6) and rewrap the whole "Remarques"
# column so every note row reads with wrapped text, matching the commit
# "ajout MCD + update journal de travail".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 5 (date 2021-11-18 / week 1) ------------------------------------
# The activity description is refined; Heures (D5, "1h30") stays as-is.
$ws.Range("C5").Value = "planification de la base de données, et création MCD"

# --- Row 6 (date 2021-11-23 / week 2) -------------------------------------
# Fill in the hours first, then the activity text, so new shared-string
# entries land in the same order the original author typed them in
# (Heures -> Activité -> Remarques), matching xl/sharedStrings.xml order.
$ws.Range("D6").Value = "3h45"
$ws.Range("C6").Value = "création de git, ajout des stories dans icescrum et fin de la création de MCD. Création du MLD"

# --- Remarques column (E) for rows 5 and 6 --------------------------------
$ws.Range("E5").Value = "On s'est planifié sur quels types de tables il y aurait, et on a également commencé a effectuer le MCD"
$ws.Range("E6").Value = "Meeting avec le professeur, nous avons crée un dépôt git et ajouté des stories dans icescrum. Le MCD est terminé. Nous avons commencé a créer le MLD"

# The whole Remarques column (E5:E28) switches to the wrap-text cell style
# (border + wrapText, same as the Activité column) so future notes wrap too.
$ws.Range("E5:E28").WrapText = $true

# Rows grow taller to fit the newly-wrapped, longer text.
$ws.Range("A5").RowHeight = 45
$ws.Range("A6").RowHeight = 60

# Last touched cell before saving, per the recorded sheet view.
$ws.Range("E8").Select()

# Best-effort: restore the recorded absolute project path shown in the
# workbook's x15ac:absPath metadata. Not all hosts expose this bit of
# OOXML bookkeeping through the Excel object model, so failures here are
# swallowed rather than aborting the rest of the edit.
try {
    $wb.Path = "C:\MA\Projet C#\github_c#\Projet-C_KGS_MMO_SGI\"
} catch {
}
